$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the two missing "unrecognized field" values to row 3
# (set E3 first so the shared-string table order matches: "MecE 265" before "MecE 260")
$ws.Range("E3").Value = "MecE 265"
$ws.Range("D3").Value = "MecE 260"

# Update the active selection to reflect the last edited cell (D3)
$ws.Range("D3").Select()
